$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 3435.0605
$ws.Range("I41").Value = 374.82352
$ws.Range("J41").Value = 6686.5625
$ws.Range("K41").Value = 374.82352
$ws.Range("L41").Value = 6686.5625
$ws.Range("M41").Value = 65.17648000000003
$ws.Range("N41").Value = -7566.5625
$ws.Range("H116").Value = 4005
$ws.Range("I116").Value = 4242.6665
$ws.Range("J116").Value = 3719.8
$ws.Range("K116").Value = 4242.6665
$ws.Range("L116").Value = 3719.8
$ws.Range("M116").Value = -800.6665000000003
$ws.Range("N116").Value = -10603.8
$ws.Range("H132").Value = 2114.2222
$ws.Range("I132").Value = 1175.7539
$ws.Range("K132").Value = 3527.2617
$ws.Range("M132").Value = -997.2617
$ws.Range("H137").Value = 1908.8983
$ws.Range("I137").Value = 1841.4318
$ws.Range("J137").Value = 2106.8
$ws.Range("K137").Value = 5524.2954
$ws.Range("L137").Value = 6320.400000000001
$ws.Range("M137").Value = -2974.2954
$ws.Range("N137").Value = -11420.4

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 7171.6665
$ws.Range("I21").Value = 5015
$ws.Range("J21").Value = 8250
$ws.Range("K21").Value = 5015
$ws.Range("L21").Value = 8250
$ws.Range("M21").Value = -4641
$ws.Range("N21").Value = -8998
$ws.Range("H26").Value = 3400
$ws.Range("I26").Value = 2211.111
$ws.Range("J26").Value = 8750
$ws.Range("K26").Value = 2211.111
$ws.Range("L26").Value = 8750
$ws.Range("M26").Value = -1881.111
$ws.Range("N26").Value = -9410
$ws.Range("H30").Value = 3288.5557
$ws.Range("I30").Value = 574.5
$ws.Range("J30").Value = 5459.8
$ws.Range("K30").Value = 574.5
$ws.Range("L30").Value = 5459.8
$ws.Range("M30").Value = -424.5
$ws.Range("N30").Value = -5759.8
$ws.Range("H74").Value = 241082.45
$ws.Range("I74").Value = 346179.34
$ws.Range("J74").Value = 61799.53
$ws.Range("K74").Value = 346179.34
$ws.Range("L74").Value = 61799.53
$ws.Range("M74").Value = -345305.34
$ws.Range("N74").Value = -63547.53
$ws.Range("H77").Value = 241082.45
$ws.Range("I77").Value = 346179.34
$ws.Range("J77").Value = 61799.53
$ws.Range("K77").Value = 1730896.7
$ws.Range("L77").Value = 308997.65
$ws.Range("M77").Value = -1726528.7
$ws.Range("N77").Value = -317733.65
$ws.Range("H113").Value = 38999.332
$ws.Range("J113").Value = 38999.332
$ws.Range("L113").Value = 38999.332
$ws.Range("N113").Value = -47677.332
$ws.Range("H122").Value = 4480.8184
$ws.Range("I122").Value = 4454.3076
$ws.Range("J122").Value = 4519.1113
$ws.Range("K122").Value = 13362.9228
$ws.Range("L122").Value = 13557.3339
$ws.Range("M122").Value = -10912.9228
$ws.Range("N122").Value = -18457.3339
$ws.Range("H132").Value = 22323.846
$ws.Range("I132").Value = 31319.943
$ws.Range("J132").Value = 3802.4707
$ws.Range("K132").Value = 93959.829
$ws.Range("L132").Value = 11407.4121
$ws.Range("M132").Value = -91429.829
$ws.Range("N132").Value = -16467.4121

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1152.3077
$ws.Range("I94").Value = 724.4
$ws.Range("J94").Value = 1735.8182
$ws.Range("K94").Value = 724.4
$ws.Range("L94").Value = 1735.8182
$ws.Range("M94").Value = -273.4
$ws.Range("N94").Value = -2637.8182

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 857.0625
$ws.Range("I16").Value = 863.3077
$ws.Range("J16").Value = 830
$ws.Range("K16").Value = 863.3077
$ws.Range("L16").Value = 830
$ws.Range("M16").Value = -576.3077
$ws.Range("N16").Value = -1404
$ws.Range("H31").Value = 2393.7576
$ws.Range("I31").Value = 1878.5238
$ws.Range("J31").Value = 3295.4167
$ws.Range("K31").Value = 1878.5238
$ws.Range("L31").Value = 3295.4167
$ws.Range("M31").Value = -1583.5238
$ws.Range("N31").Value = -3885.4167
$ws.Range("H34").Value = 2393.7576
$ws.Range("I34").Value = 1878.5238
$ws.Range("J34").Value = 3295.4167
$ws.Range("K34").Value = 1878.5238
$ws.Range("L34").Value = 3295.4167
$ws.Range("M34").Value = -1676.5238
$ws.Range("N34").Value = -3699.4167
$ws.Range("H86").Value = 4223.4614
$ws.Range("I86").Value = 1588.125
$ws.Range("K86").Value = 1588.125
$ws.Range("M86").Value = -465.125
$ws.Range("H89").Value = 4223.4614
$ws.Range("I89").Value = 1588.125
$ws.Range("K89").Value = 7940.625
$ws.Range("M89").Value = -2324.625
$ws.Range("H99").Value = 34314.742
$ws.Range("I99").Value = 60649.883
$ws.Range("J99").Value = 2336.3572
$ws.Range("K99").Value = 60649.883
$ws.Range("L99").Value = 2336.3572
$ws.Range("M99").Value = -59151.883
$ws.Range("N99").Value = -5332.3572
$ws.Range("H113").Value = 857.0625
$ws.Range("I113").Value = 863.3077
$ws.Range("J113").Value = 830
$ws.Range("K113").Value = 863.3077
$ws.Range("L113").Value = 830
$ws.Range("M113").Value = 1306.6923
$ws.Range("N113").Value = -5170
$ws.Range("H122").Value = 978
$ws.Range("I122").Value = 777
$ws.Range("J122").Value = 1128.75
$ws.Range("K122").Value = 2331
$ws.Range("L122").Value = 3386.25
$ws.Range("M122").Value = 119
$ws.Range("N122").Value = -8286.25
$ws.Range("H126").Value = 34314.742
$ws.Range("I126").Value = 60649.883
$ws.Range("J126").Value = 2336.3572
$ws.Range("K126").Value = 181949.649
$ws.Range("L126").Value = 7009.071599999999
$ws.Range("M126").Value = -179479.649
$ws.Range("N126").Value = -11949.0716
$ws.Range("H132").Value = 2024.8
$ws.Range("I132").Value = 1072.1333
$ws.Range("K132").Value = 3216.3999
$ws.Range("M132").Value = -686.3998999999999
$ws.Range("H134").Value = 2002.2
$ws.Range("I134").Value = 1198
$ws.Range("K134").Value = 3594
$ws.Range("M134").Value = -1059

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 20000196
$ws.Range("I92").Value = 20833530
$ws.Range("K92").Value = 62500590
$ws.Range("M92").Value = -62499342
$ws.Range("H122").Value = 11905302
$ws.Range("I122").Value = 17241622
$ws.Range("J122").Value = 1204.8462
$ws.Range("K122").Value = 155174598
$ws.Range("L122").Value = 10843.6158
$ws.Range("M122").Value = -155172148
$ws.Range("N122").Value = -15743.6158
$ws.Range("H132").Value = 3848.9688
$ws.Range("I132").Value = 2200.875
$ws.Range("J132").Value = 5497.0625
$ws.Range("K132").Value = 19807.875
$ws.Range("L132").Value = 49473.5625
$ws.Range("M132").Value = -17277.875
$ws.Range("N132").Value = -54533.5625

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 16693304
$ws.Range("I18").Value = 25004952
$ws.Range("J18").Value = 70006
$ws.Range("K18").Value = 25004952
$ws.Range("L18").Value = 70006
$ws.Range("M18").Value = -25004659
$ws.Range("N18").Value = -70592
$ws.Range("H29").Value = 5448.5
$ws.Range("I29").Value = 5448.5
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 5448.5
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -5158.5
$ws.Range("N29").Value = ""
$ws.Range("H102").Value = 3518.8928
$ws.Range("I102").Value = 2026.2941
$ws.Range("J102").Value = 5825.636
$ws.Range("K102").Value = 2026.2941
$ws.Range("L102").Value = 5825.636
$ws.Range("M102").Value = -404.2941000000001
$ws.Range("N102").Value = -9069.636
$ws.Range("H113").Value = 1393.5
$ws.Range("I113").Value = 1105.0834
$ws.Range("J113").Value = 1681.9166
$ws.Range("K113").Value = 1105.0834
$ws.Range("L113").Value = 1681.9166
$ws.Range("M113").Value = 1064.9166
$ws.Range("N113").Value = -6021.9166
$ws.Range("H132").Value = 3492.8298
$ws.Range("I132").Value = 3333.3794
$ws.Range("J132").Value = 3749.7222
$ws.Range("K132").Value = 10000.1382
$ws.Range("L132").Value = 11249.1666
$ws.Range("M132").Value = -7470.138199999999
$ws.Range("N132").Value = -16309.1666

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1636.7273
$ws.Range("I61").Value = 1525.5
$ws.Range("J61").Value = 1933.3334
$ws.Range("K61").Value = 1525.5
$ws.Range("L61").Value = 1933.3334
$ws.Range("M61").Value = -1323.5
$ws.Range("N61").Value = -2337.3334
$ws.Range("H113").Value = 1636.7273
$ws.Range("I113").Value = 1525.5
$ws.Range("J113").Value = 1933.3334
$ws.Range("K113").Value = 1525.5
$ws.Range("L113").Value = 1933.3334
$ws.Range("M113").Value = 644.5
$ws.Range("N113").Value = -6273.3334
$ws.Range("H132").Value = 10704.308
$ws.Range("I132").Value = 4256.5713
$ws.Range("J132").Value = 13079.789
$ws.Range("K132").Value = 12769.7139
$ws.Range("L132").Value = 39239.367
$ws.Range("M132").Value = -10239.7139
$ws.Range("N132").Value = -44299.367

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 325
$ws.Range("I113").Value = 342
$ws.Range("J113").Value = 308
$ws.Range("K113").Value = 1026
$ws.Range("L113").Value = 924
$ws.Range("M113").Value = 1144
$ws.Range("N113").Value = -5264
$ws.Range("H122").Value = 23811002
$ws.Range("I122").Value = 37038450
$ws.Range("J122").Value = 1598
$ws.Range("K122").Value = 111115350
$ws.Range("L122").Value = 4794
$ws.Range("M122").Value = -111112900
$ws.Range("N122").Value = -9694
$ws.Range("H132").Value = 1709.0426
$ws.Range("I132").Value = 1035.5946
$ws.Range("J132").Value = 4200.8
$ws.Range("K132").Value = 3106.7838
$ws.Range("L132").Value = 12602.4
$ws.Range("M132").Value = -576.7837999999997
$ws.Range("N132").Value = -17662.4
